# Regenerate merged AHB files
#
# 1) Rename the header row: the "_old" / "_new" column-name suffixes used by
#    the diff/merge tooling become version tags "_FV2310" / "_FV2404".
# 2) Turn the data range A1:U71 into a native Excel Table ("Table1") so the
#    header row gets AutoFilter buttons and the column names are tracked by
#    the table definition.
# 3) Freeze the header row (split/freeze at row 2, i.e. below row 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename header cells -------------------------------------------------
$ws.Range("A1").Value = "Segmentname_FV2310"
$ws.Range("B1").Value = "Segmentgruppe_FV2310"
$ws.Range("C1").Value = "Segment_FV2310"
$ws.Range("D1").Value = "Datenelement_FV2310"
$ws.Range("E1").Value = "Segment ID_FV2310"
$ws.Range("F1").Value = "Code_FV2310"
$ws.Range("G1").Value = "Qualifier_FV2310"
$ws.Range("H1").Value = "Beschreibung_FV2310"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2310"
$ws.Range("J1").Value = "Bedingung_FV2310"

$ws.Range("K1").Value = "diff"

$ws.Range("L1").Value = "Segmentname_FV2404"
$ws.Range("M1").Value = "Segmentgruppe_FV2404"
$ws.Range("N1").Value = "Segment_FV2404"
$ws.Range("O1").Value = "Datenelement_FV2404"
$ws.Range("P1").Value = "Segment ID_FV2404"
$ws.Range("Q1").Value = "Code_FV2404"
$ws.Range("R1").Value = "Qualifier_FV2404"
$ws.Range("S1").Value = "Beschreibung_FV2404"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2404"
$ws.Range("U1").Value = "Bedingung_FV2404"

# --- 2) Turn A1:U71 into a table named "Table1" -----------------------------
$tableRange = $ws.Range("A1:U71")
$lo = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $tableRange,
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$lo.Name = "Table1"
$lo.TableStyle = ""
$lo.ShowTableStyleRowStripes = $true
$lo.ShowTableStyleFirstColumn = $false
$lo.ShowTableStyleLastColumn = $false
$lo.ShowTableStyleColumnStripes = $false

# --- 3) Freeze panes above row 2 (keep header row visible while scrolling) --
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Return the selection to the natural A1 default once the freeze is in place.
$ws.Range("A1").Select()
